$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = -21.0933
$ws.Range("A10").Value = -20.52699999999997
$ws.Range("A12").Value = -22.45120000000002
$ws.Range("C13").Value = -12.43639999999999
$ws.Range("A18").Value = -22.30370000000002
$ws.Range("A25").Value = -22.36190000000004
